# Weekly update: insert a new data row for "Poroto granado" (Feria Lagunitas de
# Puerto Montt) right after the existing row 47, shifting the remaining
# historical rows (old 48-71) down by one to become rows 49-72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; this pushes old rows 48..71 down to 49..72.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = 45001
$ws.Range("E48").Value = 10
$ws.Range("F48").Value = 100112030
$ws.Range("G48").Value = "Poroto granado"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 50
$ws.Range("K48").Value = 40000
$ws.Range("L48").Value = 40000
$ws.Range("M48").Value = 40000
$ws.Range("N48").Value = "$/saco 25 kilos"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 1600
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
